$p = $ppt.ActivePresentation

# --- Edit 1: Slide 2, "Calibration of STICS" content placeholder ---
# Paragraph: "I could then try to calibrate STICS with parameters derived from
# FSPM simulation outputs (i.e. inverse modeling, Gaudio et al. 2021), mostly
# for the ones impossible to measure"
# Change the citation year "2021" -> "2022" (as a distinct, freshly typed
# run), keeping the rest of the sentence intact.
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$paras2 = $tr2.Paragraphs()
for ($i = 1; $i -le $paras2.Count; $i++) {
    $para = $tr2.Paragraphs($i, 1)
    if ($para.Text.IndexOf("2021), mostly for the ones impossible to measure") -ge 0) {
        $offset = $para.Text.IndexOf("2021), ")
        $target = $para.Characters($offset + 1, 7)
        $target.Text = "2022), "
    }
}

# --- Edit 2: Slide 3, "What I'd need to do" content placeholder ---
# Paragraph: "Learn how to use STICS (in January with Remi), calibrate it,
# run simulations, compare" was split across two runs ("...simulations" and
# ", compare"); merge it back into a single run.
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange
$paras3 = $tr3.Paragraphs()
for ($i = 1; $i -le $paras3.Count; $i++) {
    $para = $tr3.Paragraphs($i, 1)
    if ($para.Text.IndexOf("Learn how to use STICS") -ge 0) {
        $fullText = $para.Text.TrimEnd("`r")
        # Force a genuine text change so the run is rewritten as a single
        # contiguous run, then restore the final text.
        $para.Text = "~~~placeholder~~~"
        $para.Text = $fullText
    }
}
